# Add new module to create output files of Conversions parameters in A2 script
# -> Update the "Emission" list (column E, rows 2-21) on the "Lists" sheet
#    with the new ordering of emission-type names.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lists")

$values = @(
    "CO2e_HFC",
    "contam_agua",
    "Health",
    "RM",
    "CONHICK",
    "RESHID",
    "CO2e_TRN",
    "CO2e_WASTE",
    "turismo_residuos",
    "DAPANI",
    "CONHAB",
    "CONTUR",
    "CONVAR",
    "CO2e_PIUP",
    "CO2e_PP",
    "CO2e_AFOLU",
    "salud_residuos",
    "CO2e_sources",
    "FERT_ORG",
    "CO2e_DE"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $values[$i]
}
